$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2, 3, 4 have their A (Id), Q (Ost), R (Nord) and AC (Publik kommentar)
# values cyclically rotated: row2 <- row3, row3 <- row4, row4 <- row2.

# Capture original values before overwriting anything.
$A2 = $ws.Range("A2").Value2
$Q2 = $ws.Range("Q2").Value2
$R2 = $ws.Range("R2").Value2
$AC2 = $ws.Range("AC2").Value2

$A3 = $ws.Range("A3").Value2
$Q3 = $ws.Range("Q3").Value2
$R3 = $ws.Range("R3").Value2
$AC3 = $ws.Range("AC3").Value2

$A4 = $ws.Range("A4").Value2
$Q4 = $ws.Range("Q4").Value2
$R4 = $ws.Range("R4").Value2
$AC4 = $ws.Range("AC4").Value2

# Row 2 becomes old row 3
$ws.Range("A2").Value2 = $A3
$ws.Range("Q2").Value2 = $Q3
$ws.Range("R2").Value2 = $R3
$ws.Range("AC2").Value2 = $AC3

# Row 3 becomes old row 4
$ws.Range("A3").Value2 = $A4
$ws.Range("Q3").Value2 = $Q4
$ws.Range("R3").Value2 = $R4
$ws.Range("AC3").Value2 = $AC4

# Row 4 becomes old row 2
$ws.Range("A4").Value2 = $A2
$ws.Range("Q4").Value2 = $Q2
$ws.Range("R4").Value2 = $R2
$ws.Range("AC4").Value2 = $AC2
